$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 7.897899999999998
$ws.Range("A3").Value = -21.8869
$ws.Range("C3").Value = -11.4861
$ws.Range("D6").Value = -7.836999999999997
$ws.Range("C12").Value = -11.52139999999999
$ws.Range("A14").Value = -21.5971
$ws.Range("D19").Value = -8.749099999999993
$ws.Range("A21").Value = -20.27969999999998
$ws.Range("A23").Value = -20.48509999999997
$ws.Range("C24").Value = -12.38639999999999
$ws.Range("D24").Value = -7.861899999999999
$ws.Range("A25").Value = -21.91179999999999
$ws.Range("B25").Value = 5.464599999999999
$ws.Range("C25").Value = -12.86949999999999
$ws.Range("A26").Value = -21.14269999999997
$ws.Range("B27").Value = 5.594800000000005
$ws.Range("A29").Value = -20.79169999999998
$ws.Range("D30").Value = -7.105200000000004
$ws.Range("B31").Value = 5.205500000000001
$ws.Range("D31").Value = -8.086699999999997
$ws.Range("D33").Value = -8.0259
$ws.Range("B39").Value = 9.651400000000001
$ws.Range("D42").Value = -9.050499999999996
$ws.Range("B48").Value = 5.223000000000001
$ws.Range("C50").Value = -13.643
$ws.Range("B51").Value = 5.219900000000001
$ws.Range("B52").Value = 4.977800000000003
$ws.Range("A53").Value = -21.45539999999999
$ws.Range("C53").Value = -10.2553
$ws.Range("B55").Value = 6.025199999999995
$ws.Range("D55").Value = -7.897399999999998
$ws.Range("B56").Value = 4.912699999999997
$ws.Range("A57").Value = -22.0919
$ws.Range("B57").Value = 5.043099999999996
$ws.Range("C57").Value = -13.60089999999998
$ws.Range("D58").Value = -8.288700000000002
$ws.Range("A59").Value = -22.31130000000001
$ws.Range("C61").Value = -13.16379999999999
$ws.Range("C63").Value = -11.3054
$ws.Range("D65").Value = -8.229900000000001
$ws.Range("A69").Value = -21.64909999999999
$ws.Range("C70").Value = -12.0246
$ws.Range("D70").Value = -8.258899999999999
$ws.Range("B73").Value = 8.466199999999995
$ws.Range("D75").Value = -8.075100000000006
$ws.Range("A79").Value = -20.60340000000001
$ws.Range("A83").Value = -22.2
$ws.Range("D83").Value = -8.224299999999996
$ws.Range("C86").Value = -13.58819999999999
$ws.Range("D86").Value = -8.192399999999994
$ws.Range("B89").Value = 5.086199999999996
$ws.Range("B90").Value = 5.579600000000002
$ws.Range("A91").Value = -21.27960000000001
$ws.Range("B92").Value = 4.934799999999995
$ws.Range("A93").Value = -20.80419999999998
$ws.Range("D96").Value = -7.2874
$ws.Range("D97").Value = -8.166399999999999
$ws.Range("C98").Value = -11.67129999999999
$ws.Range("C100").Value = -12.57629999999999
$ws.Range("C102").Value = -13.5246
